# Update countries & provincias Spain
# Applies the 23-May-2020 01:05 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last refreshed" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 01:05"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 1644062
$ws.Range("C4").Value = 23165
$ws.Range("D4").Value = 397025
$ws.Range("E4").Value = 1149429
$ws.Range("G4").Value = 1254
$ws.Range("H4").Value = 97608

# --- Row 5: Brasil ---
$ws.Range("D5").Value = 135430
$ws.Range("E5").Value = 174412

# --- Row 50: Panama ---
$ws.Range("B50").Value = 10649
$ws.Range("C50").Value = 718
$ws.Range("D50").Value = 3062
$ws.Range("E50").Value = 7154
$ws.Range("G50").Value = 17
$ws.Range("H50").Value = 433

# --- Row 51: Argentina ---
$ws.Range("B51").Value = 10267
$ws.Range("C51").Value = 151
$ws.Range("D51").Value = 6275
$ws.Range("E51").Value = 3697
$ws.Range("G51").Value = 4
$ws.Range("H51").Value = 295

# --- Row 53: Chequia ---
$ws.Range("B53").Value = 8813
$ws.Range("C53").Value = 59
$ws.Range("D53").Value = 6025
$ws.Range("E53").Value = 2476

# --- Row 55: Noruega ---
$ws.Range("B55").Value = 8332
$ws.Range("C55").Value = 23
$ws.Range("E55").Value = 370

# --- Row 59: Malasia ---
$ws.Range("B59").Value = 7261
$ws.Range("C59").Value = 245
$ws.Range("D59").Value = 2007
$ws.Range("E59").Value = 5033
$ws.Range("G59").Value = 10
$ws.Range("H59").Value = 221

# --- Row 60: Australia ---
$ws.Range("B60").Value = 7137
$ws.Range("C60").Value = 78
$ws.Range("D60").Value = 5859
$ws.Range("E60").Value = 1163
$ws.Range("H60").Value = 115

# --- Row 61: Nigeria ---
$ws.Range("B61").Value = 7095
$ws.Range("C61").Value = 14
$ws.Range("D61").Value = 6479
$ws.Range("E61").Value = 515
$ws.Range("G61").Value = 1
$ws.Range("H61").Value = 101

# --- Row 75: Guinea ---
$ws.Range("E75").Value = 1473
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 19

# --- Row 120: Uruguay ---
$ws.Range("B120").Value = 753
$ws.Range("C120").Value = 4
$ws.Range("D120").Value = 603
$ws.Range("E120").Value = 130
